$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Lot Details": fill in the monthly acquisition lots for 2025 (rows
# 3-14), plus a few trailing blank-but-styled rows (15-17) that match the
# look of the existing blank rows 3-6 in the original sheet.
# ---------------------------------------------------------------------------
$lot = $wb.Worksheets.Item("Lot Details")

$dates  = @(45658,45689,45717,45748,45778,45809,45839,45870,45901,45931,45962,45992)
$prices = @(50,500,52,53,54,55,56,57,58,59,60,61)

for ($i = 0; $i -lt 12; $i++) {
    $r = $i + 3

    $lot.Cells.Item($r,1).Value = $dates[$i]
    $lot.Cells.Item($r,2).Value = $prices[$i]
    $lot.Cells.Item($r,3).Formula = "=2500"
    $lot.Cells.Item($r,4).Value = 1

    # Column A needs the same style as the other acquisition-date cells
    # (A2) -- new rows otherwise inherit the plain column style.
    $lot.Cells.Item(2,1).Copy()
    $lot.Cells.Item($r,1).PasteSpecial(-4122)

    # Column C (Cost) always carries the "bordered" numeric style used by
    # C2, regardless of row parity.
    $lot.Cells.Item(2,3).Copy()
    $lot.Cells.Item($r,3).PasteSpecial(-4122)

    # Column D (Exchange rate) alternates between the plain style (odd
    # rows, the column default -- already applied) and the "bordered"
    # style carried by D2 (even rows).
    if ($r % 2 -eq 0) {
        $lot.Cells.Item(2,4).Copy()
        $lot.Cells.Item($r,4).PasteSpecial(-4122)
    }
}

# Trailing styled-but-empty rows, matching the blank A/E cells already
# present on rows 3-6 of the original sheet.
$lot.Cells.Item(2,1).Copy()
$lot.Cells.Item(15,1).PasteSpecial(-4122)
$lot.Cells.Item(2,3).Copy()
$lot.Cells.Item(15,3).PasteSpecial(-4122)

$lot.Cells.Item(2,1).Copy()
$lot.Cells.Item(16,1).PasteSpecial(-4122)
$lot.Cells.Item(2,3).Copy()
$lot.Cells.Item(16,3).PasteSpecial(-4122)
$lot.Cells.Item(2,4).Copy()
$lot.Cells.Item(16,4).PasteSpecial(-4122)

$lot.Cells.Item(2,1).Copy()
$lot.Cells.Item(17,1).PasteSpecial(-4122)
$lot.Cells.Item(2,3).Copy()
$lot.Cells.Item(17,3).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "PFIC Details": add a new "PFIC Share Class" column, tweak the
# address text, and widen the new column.
# ---------------------------------------------------------------------------
$pfic = $wb.Worksheets.Item("PFIC Details")

$pfic.Cells.Item(1,4).Value = "PFIC Share Class"
$pfic.Cells.Item(2,2).Value = "70 Sir John Rogerson’s Quay, Dublin Ireland"
$pfic.Cells.Item(2,4).Value = "UCITS ETF (USD) Acc."

$pfic.Columns.Item(4).ColumnWidth = 26.27

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the final UI state: the
# Lot Details sheet is no longer the active tab, PFIC Details is.
# ---------------------------------------------------------------------------
$lot.Range("D22").Select() | Out-Null

$eoy = $wb.Worksheets.Item("EOY Details")
$eoy.Range("E31").Select() | Out-Null

$pfic.Activate() | Out-Null
$pfic.Range("D2").Select() | Out-Null
